$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6053187251091003
$ws.Range("B1").Value = 1.669934034347534
$ws.Range("C1").Value = 3.741125106811523
$ws.Range("D1").Value = 1.400201439857483
$ws.Range("E1").Value = 0.7307687997817993
